$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 92, shifting existing rows 92-93 down to 93-94.
# (Matches the surrounding rows' formatting, e.g. the date style on column D.)
$ws.Rows.Item(92).Insert()

# Populate the newly inserted row 92 with the new weekly record.
$ws.Range("A92").Value = 8
$ws.Range("B92").Value = "Terminal La Palmera de La Serena"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 44509
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 100112040
$ws.Range("G92").Value = "Cilantro"
$ws.Range("H92").Value = "Sin especificar"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 3200
$ws.Range("K92").Value = 1300
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 1400
$ws.Range("N92").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O92").Value = "Provincia del Elquí"
$ws.Range("P92").Value = 933
$ws.Range("Q92").Value = 1.5
$ws.Range("R92").Value = "Hortaliza"
